$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Addr='D2'; Text='29.366.94'},
    @{Addr='E2'; Text='  -0.48%  '},
    @{Addr='D3'; Text='1.846.80'},
    @{Addr='E3'; Text='  -0.25%  '},
    @{Addr='D4'; Text='0.9988'},
    @{Addr='E4'; Text='  +0.04%  '},
    @{Addr='D5'; Text='240.33'},
    @{Addr='E5'; Text='  -0.61%  '},
    @{Addr='D6'; Text='0.6311'},
    @{Addr='E6'; Text='  +0.14%  '},
    @{Addr='B8'; Text='OKB'},
    @{Addr='C8'; Text='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'},
    @{Addr='D8'; Text='47.81'},
    @{Addr='E8'; Text='  -0.47%  '},
    @{Addr='B9'; Text='Dogecoin'},
    @{Addr='C9'; Text='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'},
    @{Addr='D9'; Text='0.07547'},
    @{Addr='E9'; Text='  +0.07%  '},
    @{Addr='B10'; Text='Cardano'},
    @{Addr='C10'; Text='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'},
    @{Addr='D10'; Text='0.2955'},
    @{Addr='E10'; Text='  -0.87%  '},
    @{Addr='B11'; Text='Solana'},
    @{Addr='C11'; Text='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'},
    @{Addr='D11'; Text='24.44'},
    @{Addr='E11'; Text='  -0.01%  '},
    @{Addr='B12'; Text='TRON'},
    @{Addr='C12'; Text='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Addr='D12'; Text='0.07705'},
    @{Addr='E12'; Text='  -0.27%  '},
    @{Addr='B13'; Text='WrappedEther'},
    @{Addr='C13'; Text='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Addr='D13'; Text='1.868.92'},
    @{Addr='E13'; Text='  +1.10%  '},
    @{Addr='B14'; Text='Polkadot'},
    @{Addr='C14'; Text='https://coinranking.com/coin/25W7FG7om+polkadot-dot'},
    @{Addr='D14'; Text='4.995'},
    @{Addr='E14'; Text='  -0.26%  '},
    @{Addr='B15'; Text='Polygon'},
    @{Addr='C15'; Text='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'},
    @{Addr='D15'; Text='0.6829'},
    @{Addr='E15'; Text='  -1.35%  '},
    @{Addr='B16'; Text='ShibaInu'},
    @{Addr='C16'; Text='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'},
    @{Addr='D16'; Text='0.00001006'},
    @{Addr='E16'; Text='  +2.87%  '},
    @{Addr='B17'; Text='Litecoin'},
    @{Addr='C17'; Text='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Addr='D17'; Text='82.88'},
    @{Addr='E17'; Text='  -0.97%  '},
    @{Addr='B18'; Text='WrappedliquidstakedEther2.0'},
    @{Addr='C18'; Text='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{Addr='D18'; Text='2.122.42'},
    @{Addr='E18'; Text='  -1.24%  '},
    @{Addr='B19'; Text='Uniswap'},
    @{Addr='C19'; Text='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'},
    @{Addr='D19'; Text='6.134'},
    @{Addr='E19'; Text='  -1.81%  '},
    @{Addr='B20'; Text='WrappedBTC'},
    @{Addr='C20'; Text='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'},
    @{Addr='D20'; Text='29.397.94'},
    @{Addr='E20'; Text='  -0.58%  '},
    @{Addr='B21'; Text='BitcoinCash'},
    @{Addr='C21'; Text='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'},
    @{Addr='D21'; Text='227.44'},
    @{Addr='E21'; Text='  -2.65%  '},
    @{Addr='B22'; Text='Avalanche'},
    @{Addr='C22'; Text='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'},
    @{Addr='D22'; Text='12.45'},
    @{Addr='E22'; Text='  -0.40%  '},
    @{Addr='B23'; Text='Dai'},
    @{Addr='C23'; Text='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Addr='D23'; Text='0.9999'},
    @{Addr='E23'; Text='  +0.06%  '},
    @{Addr='B24'; Text='Chainlink'},
    @{Addr='C24'; Text='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'},
    @{Addr='D24'; Text='7.545'},
    @{Addr='E24'; Text='  -1.54%  '},
    @{Addr='B25'; Text='BinanceUSD'},
    @{Addr='C25'; Text='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'},
    @{Addr='D25'; Text='1.000'},
    @{Addr='E25'; Text='  +0.09%  '},
    @{Addr='B26'; Text='Monero'},
    @{Addr='C26'; Text='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Addr='D26'; Text='157.33'},
    @{Addr='E26'; Text='  +1.74%  '},
    @{Addr='B27'; Text='Stellar'},
    @{Addr='C27'; Text='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Addr='D27'; Text='0.1394'},
    @{Addr='E27'; Text='  +0.11%  '},
    @{Addr='B28'; Text='Cosmos'},
    @{Addr='C28'; Text='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Addr='D28'; Text='8.354'},
    @{Addr='E28'; Text='  -1.27%  '},
    @{Addr='B29'; Text='EthereumClassic'},
    @{Addr='C29'; Text='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{Addr='D29'; Text='17.64'},
    @{Addr='E29'; Text='  -0.38%  '},
    @{Addr='B30'; Text='PancakeSwap'},
    @{Addr='C30'; Text='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{Addr='D30'; Text='1.468'},
    @{Addr='E30'; Text='  -0.86%  '},
    @{Addr='D31'; Text='1.260'},
    @{Addr='E31'; Text='  +0.76%  '},
    @{Addr='B32'; Text='Hedera'},
    @{Addr='C32'; Text='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Addr='D32'; Text='0.05676'},
    @{Addr='E32'; Text='  -4.05%  '},
    @{Addr='B33'; Text='Filecoin'},
    @{Addr='C33'; Text='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Addr='D33'; Text='4.119'},
    @{Addr='E33'; Text='  +0.26%  '},
    @{Addr='B34'; Text='InternetComputer(DFINITY)'},
    @{Addr='C34'; Text='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'},
    @{Addr='D34'; Text='4.017'},
    @{Addr='E34'; Text='  -0.50%  '},
    @{Addr='B35'; Text='LidoDAOToken'},
    @{Addr='C35'; Text='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'},
    @{Addr='D35'; Text='1.844'},
    @{Addr='E35'; Text='  -1.78%  '},
    @{Addr='B36'; Text='ARBITRUM'},
    @{Addr='C36'; Text='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{Addr='D36'; Text='1.154'},
    @{Addr='E36'; Text='  -1.25%  '},
    @{Addr='B37'; Text='ImmutableX'},
    @{Addr='C37'; Text='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Addr='D37'; Text='0.7115'},
    @{Addr='E37'; Text='  -1.26%  '},
    @{Addr='B38'; Text='HuobiToken'},
    @{Addr='C38'; Text='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'},
    @{Addr='D38'; Text='2.592'},
    @{Addr='E38'; Text='  +0.14%  '},
    @{Addr='B39'; Text='Maker'},
    @{Addr='C39'; Text='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Addr='D39'; Text='1.263.47'},
    @{Addr='E39'; Text='  +1.68%  '},
    @{Addr='B40'; Text='VeChain'},
    @{Addr='C40'; Text='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Addr='D40'; Text='0.01814'},
    @{Addr='E40'; Text='  +0.74%  '},
    @{Addr='B41'; Text='MXToken'},
    @{Addr='C41'; Text='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Addr='D41'; Text='2.777'},
    @{Addr='E41'; Text='  -0.69%  '},
    @{Addr='D42'; Text='0.9098'},
    @{Addr='E42'; Text='  +0.23%  '},
    @{Addr='B43'; Text='FraxShare'},
    @{Addr='C43'; Text='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{Addr='D43'; Text='6.200'},
    @{Addr='E43'; Text='  +0.78%  '},
    @{Addr='B44'; Text='PaxDollar'},
    @{Addr='C44'; Text='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'},
    @{Addr='D44'; Text='0.9999'},
    @{Addr='E44'; Text='  +0.06%  '},
    @{Addr='B45'; Text='Quant'},
    @{Addr='C45'; Text='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'},
    @{Addr='D45'; Text='101.40'},
    @{Addr='E45'; Text='  -0.36%  '},
    @{Addr='B46'; Text='Aave'},
    @{Addr='C46'; Text='https://coinranking.com/coin/ixgUfzmLR+aave-aave'},
    @{Addr='D46'; Text='66.24'},
    @{Addr='E46'; Text='  -1.53%  '},
    @{Addr='B47'; Text='Aptos'},
    @{Addr='C47'; Text='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Addr='D47'; Text='7.068'},
    @{Addr='E47'; Text='  -4.74%  '},
    @{Addr='B48'; Text='TheSandbox'},
    @{Addr='C48'; Text='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Addr='D48'; Text='0.4044'},
    @{Addr='E48'; Text='  -0.17%  '},
    @{Addr='D49'; Text='9.104'},
    @{Addr='E49'; Text='  -0.62%  '},
    @{Addr='D50'; Text='1.682'},
    @{Addr='E50'; Text='  -1.07%  '},
    @{Addr='D51'; Text='0.1122'},
    @{Addr='E51'; Text='  +0.22%  '}
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
    $cell.ClearFormats()
}
